# Insert a new data row (new weekly record) right after the existing
# row 110, shifting all subsequent rows (old 111..170) down by one
# (new 112..171), and populate the newly inserted row 111 with the
# new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111; this shifts rows 111-170 down to 112-171
# and Excel carries the row-above formatting (incl. the date number format
# on column D) into the freshly inserted row automatically.
$ws.Rows.Item(111).Insert()

# Populate the new row 111 with the new record.
$ws.Cells.Item(111, 1).Value  = 9
$ws.Cells.Item(111, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(111, 3).Value  = "Metropolitana"
$ws.Cells.Item(111, 4).Value  = 44529
$ws.Cells.Item(111, 5).Value  = 13
$ws.Cells.Item(111, 6).Value  = 100112030
$ws.Cells.Item(111, 7).Value  = "Poroto granado"
$ws.Cells.Item(111, 8).Value  = "Sin especificar"
$ws.Cells.Item(111, 9).Value  = "Primera"
$ws.Cells.Item(111, 10).Value = 25
$ws.Cells.Item(111, 11).Value = 46000
$ws.Cells.Item(111, 12).Value = 48000
$ws.Cells.Item(111, 13).Value = 46960
$ws.Cells.Item(111, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(111, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(111, 16).Value = 1878
$ws.Cells.Item(111, 17).Value = 25
$ws.Cells.Item(111, 18).Value = "Hortaliza"
